$d = $word.ActiveDocument

# --- Locate the target paragraph / anchor text ---------------------------
$anchor = "Implementando a UML para a modelagem de requisitos"

$find = $d.Content
$ok = $find.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "anchor text not found"
}

# $find is now collapsed to the found anchor text; Start/End bound it.
$p0 = $find.End

# The two characters right after the anchor are ", " (comma + space) which
# precede the old "utiliza-se..." continuation. We replace that ", " with
# the newly authored sentences, and keep the "u" of "utiliza-se" merged
# into the tail of our insertion (" U" + "tiliza-se...").
$commaSpace = $d.Range($p0, $p0 + 2)
if ($commaSpace.Text -ne ", ") {
    throw ("unexpected text after anchor: '" + $commaSpace.Text + "'")
}
$commaSpace.Delete()

# --- Helper: insert text at a collapsed point and stamp explicit ----------
# character formatting on exactly the inserted span, so it forms its own
# run instead of merging into whatever precedes it.
function Insert-Formatted([int]$pos, [string]$text, [bool]$small) {
    $ip = $d.Range($pos, $pos)
    $ip.InsertAfter($text)
    $newEnd = $pos + $text.Length
    $r = $d.Range($pos, $newEnd)
    if ($small) {
        $r.Font.Size = 10
        $r.Font.SizeBi = 10
    } else {
        $r.Font.Size = 11
        $r.Font.SizeBi = 11
    }
    return $newEnd
}

$pos = $p0

# 1) ", de acordo com Grady Booch (um dos criadores da UML)" - sz20
$seg1 = ", de acordo com Grady Booch (um dos criadores da UML)"
$pos = Insert-Formatted $pos $seg1 $true

# 2) ", ele " - default (no explicit sz override)
$seg2 = ", ele "
$pos = Insert-Formatted $pos $seg2 $false

# 3) "enxerga a UML ... incluindo os requisitos. U" - sz20
$seg3 = "enxerga a UML como uma linguagem padronizada para visualizar, especificar, construir e documentar os artefatos de um sistema de software " + [char]0x2014 + " incluindo os requisitos. U"
$pos = Insert-Formatted $pos $seg3 $true

# The remaining original text ("tiliza-se quatro tipos ... completa do sistema.")
# is untouched and keeps its original run/formatting.

Write-Output $d.Paragraphs(3).Range.Text
